$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Beta) values ---
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 40.17361422455784
$ws.Range("G2").Value = 3.023936917708523
$ws.Range("H2").Value = 113.2997885567166
$ws.Range("I2").Value = 4.26388169043677
$ws.Range("J2").Value = 4.262498394416682
$ws.Range("K2").Value = 4.265849736929503
$ws.Range("L2").Value = 0.2394993062799019
$ws.Range("M2").Value = 0.2324654891153624
$ws.Range("N2").Value = 0.2533026716387036

# --- Update row 3 (Gamma) values ---
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.2496598256557464
$ws.Range("G3").Value = 0.2494293762464635
$ws.Range("H3").Value = 0.2498902723954897
$ws.Range("I3").Value = 0.2235333391462987
$ws.Range("J3").Value = 0.2233209996341272
$ws.Range("K3").Value = 0.2237438655223813
$ws.Range("L3").Value = 0.2472455819478694
$ws.Range("M3").Value = 0.2470166395069133
$ws.Range("N3").Value = 0.2474746485431898

# --- Add new row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 40.42327405021359
$ws.Range("G4").Value = 3.273366293954987
$ws.Range("H4").Value = 113.549678829112
$ws.Range("I4").Value = 4.487415029583068
$ws.Range("J4").Value = 4.48581939405081
$ws.Range("K4").Value = 4.489593602451884
$ws.Range("L4").Value = 0.4867448882277714
$ws.Range("M4").Value = 0.4794821286222757
$ws.Range("N4").Value = 0.5007773201818935

# Copy the style/formatting from A3 (which carries the bordered/centered
# header-like style used for the index column) onto the new A4 cell so it
# matches the other rows in that column.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
